$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1155.4783
$ws.Range("J17").Value = 1203.9546
$ws.Range("L17").Value = 3611.8638
$ws.Range("N17").Value = -3947.8638

$ws.Range("H32").Value = 1633.4166
$ws.Range("I32").Value = 1200.4
$ws.Range("J32").Value = 1942.7142
$ws.Range("K32").Value = 1200.4
$ws.Range("L32").Value = 1942.7142
$ws.Range("M32").Value = -874.4000000000001
$ws.Range("N32").Value = -2594.7142

$ws.Range("H33").Value = 285.8
$ws.Range("I33").Value = 281.08334
$ws.Range("J33").Value = 304.66666
$ws.Range("K33").Value = 281.08334
$ws.Range("L33").Value = 304.66666
$ws.Range("M33").Value = -52.08334000000002
$ws.Range("N33").Value = -762.66666

$ws.Range("H43").Value = 7937350.5
$ws.Range("I43").Value = 874.75
$ws.Range("K43").Value = 874.75
$ws.Range("M43").Value = -805.75

$ws.Range("H53").Value = 1718.125
$ws.Range("I53").Value = 2660
$ws.Range("J53").Value = 148.33333
$ws.Range("K53").Value = 2660
$ws.Range("L53").Value = 148.33333
$ws.Range("M53").Value = -2023
$ws.Range("N53").Value = -1422.33333

$ws.Range("H58").Value = 3896.111
$ws.Range("I58").Value = 688.3333
$ws.Range("J58").Value = 5500
$ws.Range("K58").Value = 2064.9999
$ws.Range("L58").Value = 16500
$ws.Range("M58").Value = -1914.9999
$ws.Range("N58").Value = -16800

$ws.Range("H112").Value = 2294.9443
$ws.Range("I112").Value = 749.8333
$ws.Range("K112").Value = 2249.4999
$ws.Range("M112").Value = -1141.4999

$ws.Range("H132").Value = 6445.913
$ws.Range("I132").Value = 3584.75
$ws.Range("J132").Value = 12985.714
$ws.Range("K132").Value = 10754.25
$ws.Range("L132").Value = 38957.142
$ws.Range("M132").Value = -8224.25
$ws.Range("N132").Value = -44017.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3843.513
$ws.Range("I32").Value = 3684.946
$ws.Range("K32").Value = 3684.946
$ws.Range("M32").Value = -3397.946

$ws.Range("H123").Value = 68333.336
$ws.Range("J123").Value = 68333.336
$ws.Range("L123").Value = 68333.336
$ws.Range("N123").Value = -78133.336

$ws.Range("H132").Value = 3025.879
$ws.Range("I132").Value = 2806.5
$ws.Range("K132").Value = 8419.5
$ws.Range("M132").Value = -5889.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1244.7084
$ws.Range("I20").Value = 1120.2222
$ws.Range("J20").Value = 1618.1666
$ws.Range("K20").Value = 1120.2222
$ws.Range("L20").Value = 1618.1666
$ws.Range("M20").Value = -873.2221999999999
$ws.Range("N20").Value = -2112.1666

$ws.Range("H105").Value = 125002170
$ws.Range("I105").Value = 142859380
$ws.Range("J105").Value = 1741
$ws.Range("K105").Value = 142859380
$ws.Range("L105").Value = 1741
$ws.Range("M105").Value = -142857633
$ws.Range("N105").Value = -5235

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H12").Value = 5000746
$ws.Range("I12").Value = 992.5
$ws.Range("J12").Value = 10000500
$ws.Range("K12").Value = 992.5
$ws.Range("L12").Value = 10000500
$ws.Range("M12").Value = -822.5
$ws.Range("N12").Value = -10000840

$ws.Range("H31").Value = 1312.975
$ws.Range("I31").Value = 1012.17645
$ws.Range("J31").Value = 1535.3043
$ws.Range("K31").Value = 1012.17645
$ws.Range("L31").Value = 1535.3043
$ws.Range("M31").Value = -717.17645
$ws.Range("N31").Value = -2125.3043

$ws.Range("H34").Value = 1312.975
$ws.Range("I34").Value = 1012.17645
$ws.Range("J34").Value = 1535.3043
$ws.Range("K34").Value = 1012.17645
$ws.Range("L34").Value = 1535.3043
$ws.Range("M34").Value = -810.17645
$ws.Range("N34").Value = -1939.3043

$ws.Range("H35").Value = 832
$ws.Range("I35").Value = 832
$ws.Range("K35").Value = 832
$ws.Range("M35").Value = -538

$ws.Range("H132").Value = 3205.2693
$ws.Range("I132").Value = 3115.2727
$ws.Range("J132").Value = 3700.25
$ws.Range("K132").Value = 9345.8181
$ws.Range("L132").Value = 11100.75
$ws.Range("M132").Value = -6815.8181
$ws.Range("N132").Value = -16160.75

$ws.Range("H134").Value = 1778.7916
$ws.Range("I134").Value = 1638.3889
$ws.Range("K134").Value = 4915.1667
$ws.Range("M134").Value = -2380.1667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20834724
$ws.Range("I131").Value = 111111910
$ws.Range("J131").Value = 1525.7949
$ws.Range("K131").Value = 333335730
$ws.Range("L131").Value = 4577.384700000001
$ws.Range("M131").Value = -333330690
$ws.Range("N131").Value = -14657.3847

$ws.Range("H139").Value = 1841.3636
$ws.Range("I139").Value = 1747.4615
$ws.Range("J139").Value = 1977
$ws.Range("K139").Value = 5242.3845
$ws.Range("L139").Value = 5931
$ws.Range("M139").Value = -102.3845000000001
$ws.Range("N139").Value = -16211

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 30005086
$ws.Range("I70").Value = 31254762
$ws.Range("J70").Value = 28576886
$ws.Range("K70").Value = 31254762
$ws.Range("L70").Value = 28576886
$ws.Range("M70").Value = -31254492
$ws.Range("N70").Value = -28577426

$ws.Range("H73").Value = 30005086
$ws.Range("I73").Value = 31254762
$ws.Range("J73").Value = 28576886
$ws.Range("K73").Value = 31254762
$ws.Range("L73").Value = 28576886
$ws.Range("M73").Value = -31253826
$ws.Range("N73").Value = -28578758

$ws.Range("H86").Value = 32658.75
$ws.Range("J86").Value = 32658.75
$ws.Range("L86").Value = 32658.75
$ws.Range("N86").Value = -35030.75

$ws.Range("H89").Value = 32658.75
$ws.Range("J89").Value = 32658.75
$ws.Range("L89").Value = 97976.25
$ws.Range("N89").Value = -109832.25

$ws.Range("H126").Value = 2161.6155
$ws.Range("I126").Value = 1757
$ws.Range("K126").Value = 5271
$ws.Range("M126").Value = -2801

$ws.Range("H132").Value = 2416.0715
$ws.Range("I132").Value = 2040.619
$ws.Range("K132").Value = 6121.857
$ws.Range("M132").Value = -3591.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 444.53333
$ws.Range("I55").Value = 76.42856999999999
$ws.Range("J55").Value = 766.625
$ws.Range("K55").Value = 76.42856999999999
$ws.Range("L55").Value = 766.625
$ws.Range("M55").Value = 96.57143000000001
$ws.Range("N55").Value = -1112.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 58464.5
$ws.Range("J123").Value = 58464.5
$ws.Range("L123").Value = 58464.5
$ws.Range("N123").Value = -68264.5
